# Update "想去人数" (want-to-go count) values in column F for the
# exhibitions sheet ("展览") and the combined "全部类型" sheet, matching
# the refreshed scrape snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item(1)   # 展览
$wsAllTypes    = $wb.Worksheets.Item(4)   # 全部类型

# 展览 sheet: rows 2-10 -> column F
$wsExhibitions.Range("F2").Value  = 177
$wsExhibitions.Range("F3").Value  = 663
$wsExhibitions.Range("F4").Value  = 27
$wsExhibitions.Range("F5").Value  = 225
$wsExhibitions.Range("F6").Value  = 1709
$wsExhibitions.Range("F7").Value  = 39
$wsExhibitions.Range("F8").Value  = 3216
$wsExhibitions.Range("F9").Value  = 456
$wsExhibitions.Range("F10").Value = 760

# 全部类型 sheet: same events, offset by one row after row 4
# (row 5 is a 演出 entry that is unaffected)
$wsAllTypes.Range("F2").Value  = 177
$wsAllTypes.Range("F3").Value  = 663
$wsAllTypes.Range("F4").Value  = 27
$wsAllTypes.Range("F6").Value  = 225
$wsAllTypes.Range("F7").Value  = 1709
$wsAllTypes.Range("F8").Value  = 39
$wsAllTypes.Range("F9").Value  = 3216
$wsAllTypes.Range("F10").Value = 456
$wsAllTypes.Range("F11").Value = 760
